# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the newly scraped counts.

$wb = $excel.ActiveWorkbook

# Map: sheet name -> hashtable of row -> new value
$updates = @{
    "展览" = @{
        4  = 1211
        10 = 3416
        14 = 38
        18 = 704
        19 = 203
        24 = 2509
        25 = 5020
        29 = 1287
        35 = 95
        38 = 455
        41 = 450
        43 = 460
    }
    "全部类型" = @{
        4  = 1211
        10 = 3416
        15 = 38
        19 = 704
        20 = 203
        25 = 2509
        26 = 5020
        30 = 1287
        36 = 95
        39 = 455
        42 = 450
        44 = 460
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}

$wb.Save()
